$d = $word.ActiveDocument

# Replacement 1: "as compared to" -> "in comparison to the relatively straightforward capture of"
$d.Content.Find.Execute("as compared to", $true, $false, $false, $false, $false, $true, 1, $false, "in comparison to the relatively straightforward capture of", 2)

# Replacement 2: fix "rather then individuals." typo and append new sentence content
$d.Content.Find.Execute("rather then individuals. ", $true, $false, $false, $false, $false, $true, 1, $false, "rather than individuals, furthermore concluded that ZIP area(postcodes in the UK) were the preferred form along with median household income or percent poverty as the measures for SES, although relative education state was also considered relevant in majority of studies . ", 2)
